$d = $word.ActiveDocument

# The requirement line ("LOB1038: ... (Requisito fraco)") used to be
# followed by three paragraphs: a blank paragraph, a "Ver no Jupiter
# Salvar em pdf Salvar em docx" line, and a "© 2020 ..." copyright line.
# Those three paragraphs (their text and their paragraph marks) are being
# removed, so the requirement line is directly followed by the blank
# paragraph that used to sit right before the trailing page-break
# paragraph.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "LOB1038: Física Experimental I (Requisito fraco)") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'LOB1038 ... (Requisito fraco)' paragraph"
}

$startPara = $target.Next()             # blank paragraph right after the requirement line
$endPara = $startPara.Next().Next()     # the "© 2020 ..." copyright paragraph

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
